{"js": "const replacements = [\n  [\"780\u00d78=\", \"676\u00d75=\"],\n  [\"521\u00d75=\", \"274\u00d76=\"],\n  [\"574\u00d76=\", \"674\u00d76=\"],\n  [\"403\u00d75=\", \"851\u00d72=\"],\n  [\"401\u00d79=\", \"908\u00d76=\"],\n  [\"763\u00d79=\", \"509\u00d72=\"],\n  [\"844\u00d79=\", \"667\u00d76=\"],\n  [\"480\u00d75=\", \"171\u00d76=\"],\n  [\"305\u00d79=\", \"338\u00d73=\"],\n  [\"700\u00d73=\", \"312\u00d73=\"],\n  [\"762\u00d78=\", \"480\u00d73=\"],\n  [\"415\u00d77=\", \"623\u00d73=\"],\n  [\"358\u00d74=\", \"268\u00d79=\"],\n  [\"671\u00d73=\", \"586\u00d77=\"],\n  [\"572\u00d77=\", \"206\u00d75=\"],\n  [\"229\u00d73=\", \"295\u00d76=\"],\n  [\"462\u00d79=\", \"461\u00d77=\"],\n  [\"114\u00d76=\", \"272\u00d72=\"],\n  [\"244\u00d78=\", \"895\u00d72=\"],\n  [\"769\u00d78=\", \"634\u00d72=\"],\n  [\"621\u00d79=\", \"846\u00d75=\"],\n  [\"436\u00d78=\", \"710\u00d73=\"],\n  [\"895\u00d79=\", \"770\u00d73=\"],\n  [\"982\u00d75=\", \"912\u00d73=\"],\n  [\"599\u00d76=\", \"725\u00d77=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Not found: \" + oldText);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"780\u00d78=\", \"676\u00d75=\"),\n    @(\"521\u00d75=\", \"274\u00d76=\"),\n    @(\"574\u00d76=\", \"674\u00d76=\"),\n    @(\"403\u00d75=\", \"851\u00d72=\"),\n    @(\"401\u00d79=\", \"908\u00d76=\"),\n    @(\"763\u00d79=\", \"509\u00d72=\"),\n    @(\"844\u00d79=\", \"667\u00d76=\"),\n    @(\"480\u00d75=\", \"171\u00d76=\"),\n    @(\"305\u00d79=\", \"338\u00d73=\"),\n    @(\"700\u00d73=\", \"312\u00d73=\"),\n    @(\"762\u00d78=\", \"480\u00d73=\"),\n    @(\"415\u00d77=\", \"623\u00d73=\"),\n    @(\"358\u00d74=\", \"268\u00d79=\"),\n    @(\"671\u00d73=\", \"586\u00d77=\"),\n    @(\"572\u00d77=\", \"206\u00d75=\"),\n    @(\"229\u00d73=\", \"295\u00d76=\"),\n    @(\"462\u00d79=\", \"461\u00d77=\"),\n    @(\"114\u00d76=\", \"272\u00d72=\"),\n    @(\"244\u00d78=\", \"895\u00d72=\"),\n    @(\"769\u00d78=\", \"634\u00d72=\"),\n    @(\"621\u00d79=\", \"846\u00d75=\"),\n    @(\"436\u00d78=\", \"710\u00d73=\"),\n    @(\"895\u00d79=\", \"770\u00d73=\"),\n    @(\"982\u00d75=\", \"912\u00d73=\"),\n    @(\"599\u00d76=\", \"725\u00d77=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
